# The diagram (SmartArt) graphic frame on slide 1 is moved further down the
# slide: its vertical offset goes from 1407600 EMU to 2847600 EMU (a delta of
# 1440000 EMU = 113.3858267716535 pt). PowerPoint's COM object model works in
# points (1 pt = 12700 EMU), so convert the target EMU offset to points.
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

$targetTopEmu = 2847600
$sh.Top = $targetTopEmu / 12700
